$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hello!"
$ws.Range("B2").Value = "world?"
